$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I, shifting I:T -> J:U
$ws.Columns("I:I").Insert()

# Resize the new "Monthly" column and narrow spacer column to match the
# finished layout (H holds the monthly figures, I is the thin spacer that
# used to be column H, J is the data-table label column that used to be I).
$ws.Columns("H:H").ColumnWidth = 10.5
$ws.Columns("I:I").ColumnWidth = 1.8333333333333333
$ws.Columns("J:J").ColumnWidth = 10

# Header for new "Monthly" column (bold, right-aligned)
$ws.Range("H3").Value = "Monthly"
$ws.Range("H3").Font.Bold = $true
$ws.Range("H3").HorizontalAlignment = -4152

# Monthly formulas (annual / 12) for each age-group row
$ws.Range("H7").Formula = "=G7/12"
$ws.Range("H9").Formula = "=G9/12"
$ws.Range("H11").Formula = "=G11/12"
$ws.Range("H13").Formula = "=G13/12"
$ws.Range("H15").Formula = "=G15/12"

# Blank spacer cell in the header row picks up the bold "Age group" header style
$ws.Range("I3").Font.Bold = $true

# Blank spacer cells that pick up the same currency number format as column G/H
$ws.Range("H5").NumberFormat = $ws.Range("G5").NumberFormat
$ws.Range("I5").NumberFormat = $ws.Range("G5").NumberFormat
$ws.Range("I7").NumberFormat = $ws.Range("G7").NumberFormat
$ws.Range("I9").NumberFormat = $ws.Range("G9").NumberFormat
$ws.Range("I11").NumberFormat = $ws.Range("G11").NumberFormat
$ws.Range("I13").NumberFormat = $ws.Range("G13").NumberFormat
$ws.Range("I15").NumberFormat = $ws.Range("G15").NumberFormat
